# Apply updates to the "Inscricoes" sheet per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Cells.Item(5, 5).Value = 87
$ws.Cells.Item(6, 5).Value = 30
$ws.Cells.Item(10, 5).Value = 283
$ws.Cells.Item(10, 6).Value = 143
$ws.Cells.Item(10, 8).Value = 143
$ws.Cells.Item(11, 5).Value = 200
$ws.Cells.Item(12, 5).Value = 294
$ws.Cells.Item(14, 5).Value = 84
$ws.Cells.Item(16, 5).Value = 121
$ws.Cells.Item(16, 6).Value = 71
$ws.Cells.Item(16, 8).Value = 71
$ws.Cells.Item(17, 5).Value = 56
$ws.Cells.Item(18, 5).Value = 41
$ws.Cells.Item(20, 5).Value = 64
$ws.Cells.Item(20, 6).Value = 26
$ws.Cells.Item(20, 8).Value = 26
$ws.Cells.Item(21, 5).Value = 89
$ws.Cells.Item(21, 6).Value = 52
$ws.Cells.Item(21, 8).Value = 52
$ws.Cells.Item(22, 5).Value = 113
$ws.Cells.Item(22, 6).Value = 67
$ws.Cells.Item(22, 8).Value = 67
$ws.Cells.Item(23, 5).Value = 117
$ws.Cells.Item(24, 5).Value = 135
$ws.Cells.Item(24, 6).Value = 71
$ws.Cells.Item(24, 8).Value = 71
$ws.Cells.Item(25, 5).Value = 142
$ws.Cells.Item(26, 5).Value = 80
$ws.Cells.Item(27, 5).Value = 188
$ws.Cells.Item(27, 6).Value = 106
$ws.Cells.Item(27, 8).Value = 106
$ws.Cells.Item(28, 5).Value = 119
$ws.Cells.Item(28, 6).Value = 43
$ws.Cells.Item(28, 8).Value = 43
$ws.Cells.Item(29, 5).Value = 119
$ws.Cells.Item(30, 5).Value = 138
$ws.Cells.Item(30, 6).Value = 81
$ws.Cells.Item(30, 8).Value = 81
$ws.Cells.Item(31, 5).Value = 49
$ws.Cells.Item(33, 5).Value = 182
$ws.Cells.Item(33, 6).Value = 98
$ws.Cells.Item(33, 8).Value = 98
$ws.Cells.Item(34, 5).Value = 137
$ws.Cells.Item(34, 6).Value = 88
$ws.Cells.Item(34, 8).Value = 88
$ws.Cells.Item(35, 6).Value = 60
$ws.Cells.Item(35, 8).Value = 60
$ws.Cells.Item(36, 5).Value = 40
$ws.Cells.Item(37, 5).Value = 101
$ws.Cells.Item(38, 5).Value = 59
$ws.Cells.Item(40, 5).Value = 175
$ws.Cells.Item(41, 5).Value = 238
$ws.Cells.Item(42, 5).Value = 224
$ws.Cells.Item(42, 6).Value = 120
$ws.Cells.Item(42, 8).Value = 120
$ws.Cells.Item(44, 5).Value = 184
$ws.Cells.Item(44, 6).Value = 105
$ws.Cells.Item(44, 8).Value = 105
$ws.Cells.Item(45, 5).Value = 71
$ws.Cells.Item(46, 5).Value = 180
$ws.Cells.Item(47, 5).Value = 280
$ws.Cells.Item(48, 5).Value = 131
$ws.Cells.Item(49, 5).Value = 161
$ws.Cells.Item(49, 6).Value = 79
$ws.Cells.Item(49, 8).Value = 79
$ws.Cells.Item(50, 5).Value = 135
$ws.Cells.Item(50, 6).Value = 56
$ws.Cells.Item(50, 8).Value = 56
$ws.Cells.Item(51, 5).Value = 131
$ws.Cells.Item(52, 5).Value = 17
$ws.Cells.Item(52, 6).Value = 6
$ws.Cells.Item(52, 8).Value = 6